$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(12, 9).Value = 'sd'
$ws.Cells.Item(12, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(13, 9).Value = 'sd'
$ws.Cells.Item(13, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(16, 9).Value = '%'
$ws.Cells.Item(16, 10).Value = 'Uninterpretable'
$ws.Cells.Item(42, 9).Value = 'ba'
$ws.Cells.Item(42, 10).Value = 'Appreciation'
$ws.Cells.Item(64, 9).Value = 'b'
$ws.Cells.Item(64, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(69, 9).Value = 'aa'
$ws.Cells.Item(69, 10).Value = 'Agree/Accept'
$ws.Cells.Item(80, 9).Value = 'sv'
$ws.Cells.Item(80, 10).Value = 'Statement-opinion'
$ws.Cells.Item(117, 9).Value = 'sd'
$ws.Cells.Item(117, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(134, 9).Value = 'b'
$ws.Cells.Item(134, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(141, 9).Value = 'aa'
$ws.Cells.Item(141, 10).Value = 'Agree/Accept'
$ws.Cells.Item(149, 9).Value = 'aa'
$ws.Cells.Item(149, 10).Value = 'Agree/Accept'
$ws.Cells.Item(153, 9).Value = 'sd'
$ws.Cells.Item(153, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(156, 9).Value = 'sv'
$ws.Cells.Item(156, 10).Value = 'Statement-opinion'
$ws.Cells.Item(164, 9).Value = 'sv'
$ws.Cells.Item(164, 10).Value = 'Statement-opinion'
$ws.Cells.Item(171, 9).Value = 'sd'
$ws.Cells.Item(171, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(174, 9).Value = 'sv'
$ws.Cells.Item(174, 10).Value = 'Statement-opinion'
$ws.Cells.Item(176, 9).Value = 'aa'
$ws.Cells.Item(176, 10).Value = 'Agree/Accept'
$ws.Cells.Item(178, 9).Value = 'sd'
$ws.Cells.Item(178, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(185, 9).Value = 'sd'
$ws.Cells.Item(185, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(213, 9).Value = 'aa'
$ws.Cells.Item(213, 10).Value = 'Agree/Accept'
$ws.Cells.Item(223, 9).Value = 'sd'
$ws.Cells.Item(223, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(230, 9).Value = '%'
$ws.Cells.Item(230, 10).Value = 'Uninterpretable'
$ws.Cells.Item(231, 9).Value = 'sd'
$ws.Cells.Item(231, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(237, 9).Value = 'sd'
$ws.Cells.Item(237, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(239, 9).Value = 'sv'
$ws.Cells.Item(239, 10).Value = 'Statement-opinion'
$ws.Cells.Item(241, 9).Value = 'ba'
$ws.Cells.Item(241, 10).Value = 'Appreciation'
$ws.Cells.Item(242, 9).Value = 'qy'
$ws.Cells.Item(242, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(243, 9).Value = 'sv'
$ws.Cells.Item(243, 10).Value = 'Statement-opinion'
$ws.Cells.Item(255, 9).Value = 'sd'
$ws.Cells.Item(255, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(260, 9).Value = '%'
$ws.Cells.Item(260, 10).Value = 'Uninterpretable'
$ws.Cells.Item(272, 9).Value = 'aa'
$ws.Cells.Item(272, 10).Value = 'Agree/Accept'
$ws.Cells.Item(290, 9).Value = 'sv'
$ws.Cells.Item(290, 10).Value = 'Statement-opinion'
$ws.Cells.Item(300, 9).Value = 'sd'
$ws.Cells.Item(300, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(308, 9).Value = 'ba'
$ws.Cells.Item(308, 10).Value = 'Appreciation'
$ws.Cells.Item(322, 9).Value = 'sd'
$ws.Cells.Item(322, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(325, 9).Value = 'sd'
$ws.Cells.Item(325, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(327, 9).Value = 'b'
$ws.Cells.Item(327, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(351, 9).Value = 'sv'
$ws.Cells.Item(351, 10).Value = 'Statement-opinion'
$ws.Cells.Item(352, 9).Value = 'sd'
$ws.Cells.Item(352, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(363, 9).Value = 'sd'
$ws.Cells.Item(363, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(370, 9).Value = 'sd'
$ws.Cells.Item(370, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(372, 9).Value = 'sd'
$ws.Cells.Item(372, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(384, 9).Value = 'sd'
$ws.Cells.Item(384, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(385, 9).Value = 'sd'
$ws.Cells.Item(385, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(388, 9).Value = '%'
$ws.Cells.Item(388, 10).Value = 'Uninterpretable'
$ws.Cells.Item(389, 9).Value = 'sd'
$ws.Cells.Item(389, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(390, 9).Value = '%'
$ws.Cells.Item(390, 10).Value = 'Uninterpretable'
$ws.Cells.Item(398, 9).Value = 'sd'
$ws.Cells.Item(398, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(400, 9).Value = 'sd'
$ws.Cells.Item(400, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(411, 9).Value = 'sv'
$ws.Cells.Item(411, 10).Value = 'Statement-opinion'
$ws.Cells.Item(425, 9).Value = 'sv'
$ws.Cells.Item(425, 10).Value = 'Statement-opinion'
$ws.Cells.Item(426, 9).Value = 'b'
$ws.Cells.Item(426, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(427, 9).Value = 'sd'
$ws.Cells.Item(427, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(428, 9).Value = 'sv'
$ws.Cells.Item(428, 10).Value = 'Statement-opinion'
$ws.Cells.Item(437, 9).Value = 'sv'
$ws.Cells.Item(437, 10).Value = 'Statement-opinion'
$ws.Cells.Item(443, 9).Value = 'ba'
$ws.Cells.Item(443, 10).Value = 'Appreciation'
$ws.Cells.Item(444, 9).Value = 'b'
$ws.Cells.Item(444, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(448, 9).Value = 'sv'
$ws.Cells.Item(448, 10).Value = 'Statement-opinion'
$ws.Cells.Item(456, 9).Value = 'sd'
$ws.Cells.Item(456, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(458, 9).Value = 'b'
$ws.Cells.Item(458, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(461, 9).Value = 'sv'
$ws.Cells.Item(461, 10).Value = 'Statement-opinion'
$ws.Cells.Item(464, 9).Value = 'sv'
$ws.Cells.Item(464, 10).Value = 'Statement-opinion'
$ws.Cells.Item(468, 9).Value = 'sv'
$ws.Cells.Item(468, 10).Value = 'Statement-opinion'
$ws.Cells.Item(481, 9).Value = 'sv'
$ws.Cells.Item(481, 10).Value = 'Statement-opinion'
$ws.Cells.Item(495, 9).Value = '%'
$ws.Cells.Item(495, 10).Value = 'Uninterpretable'
$ws.Cells.Item(503, 9).Value = 'sd'
$ws.Cells.Item(503, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(507, 9).Value = 'sd'
$ws.Cells.Item(507, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(518, 9).Value = 'sd'
$ws.Cells.Item(518, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(537, 9).Value = 'sd'
$ws.Cells.Item(537, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(540, 9).Value = 'b'
$ws.Cells.Item(540, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(541, 9).Value = 'aa'
$ws.Cells.Item(541, 10).Value = 'Agree/Accept'
$ws.Cells.Item(552, 9).Value = 'ba'
$ws.Cells.Item(552, 10).Value = 'Appreciation'
$ws.Cells.Item(555, 9).Value = 'sv'
$ws.Cells.Item(555, 10).Value = 'Statement-opinion'
$ws.Cells.Item(579, 9).Value = '%'
$ws.Cells.Item(579, 10).Value = 'Uninterpretable'
$ws.Cells.Item(586, 9).Value = 'sv'
$ws.Cells.Item(586, 10).Value = 'Statement-opinion'
$ws.Cells.Item(591, 9).Value = 'ba'
$ws.Cells.Item(591, 10).Value = 'Appreciation'
$ws.Cells.Item(600, 9).Value = 'sv'
$ws.Cells.Item(600, 10).Value = 'Statement-opinion'
$ws.Cells.Item(617, 9).Value = 'aa'
$ws.Cells.Item(617, 10).Value = 'Agree/Accept'
$ws.Cells.Item(619, 9).Value = 'b'
$ws.Cells.Item(619, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(630, 9).Value = 'sd'
$ws.Cells.Item(630, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(647, 9).Value = 'sv'
$ws.Cells.Item(647, 10).Value = 'Statement-opinion'
$ws.Cells.Item(648, 9).Value = 'sd'
$ws.Cells.Item(648, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(649, 9).Value = 'sv'
$ws.Cells.Item(649, 10).Value = 'Statement-opinion'
$ws.Cells.Item(672, 9).Value = 'sv'
$ws.Cells.Item(672, 10).Value = 'Statement-opinion'
$ws.Cells.Item(673, 9).Value = 'ba'
$ws.Cells.Item(673, 10).Value = 'Appreciation'
$ws.Cells.Item(675, 9).Value = 'sd'
$ws.Cells.Item(675, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(706, 9).Value = 'sd'
$ws.Cells.Item(706, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(708, 9).Value = 'aa'
$ws.Cells.Item(708, 10).Value = 'Agree/Accept'
$ws.Cells.Item(711, 9).Value = 'sv'
$ws.Cells.Item(711, 10).Value = 'Statement-opinion'
$ws.Cells.Item(719, 9).Value = 'sd'
$ws.Cells.Item(719, 10).Value = 'Statement-non-opinion'
